$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix header typos: "Precious" -> "Precise", "Genere" -> "Genre"
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "Precise" | Out-Null
$ws.Range("G1").Value = "Genre" | Out-Null

# ---------------------------------------------------------------------
# 2. Re-sort dataset: swap the song/bpm/predicted-bpm/genre fields
#    between row 2 (ID 0) and row 4 (ID 2) -- the ID/error/bool columns
#    stay put, only Song/BPM/BPM_Predict/Genre move.
# ---------------------------------------------------------------------
$b2 = $ws.Range("B2").Value2
$c2 = $ws.Range("C2").Value2
$d2 = $ws.Range("D2").Value2
$g2 = $ws.Range("G2").Value2

$b4 = $ws.Range("B4").Value2
$c4 = $ws.Range("C4").Value2
$d4 = $ws.Range("D4").Value2
$g4 = $ws.Range("G4").Value2

$ws.Range("B2").Value = $b4 | Out-Null
$ws.Range("C2").Value = $c4 | Out-Null
$ws.Range("D2").Value = $d4 | Out-Null
$ws.Range("G2").Value = $g4 | Out-Null

$ws.Range("B4").Value = $b2 | Out-Null
$ws.Range("C4").Value = $c2 | Out-Null
$ws.Range("D4").Value = $d2 | Out-Null
$ws.Range("G4").Value = $g2 | Out-Null

# ---------------------------------------------------------------------
# 3. Unify formatting across the data rows (2-17):
#    - Column B (song title) takes on the "Microsoft YaHei" font that
#      rows 3/4/12 already used.
#    - Column F (the boolean flag) gets the built-in "Good" cell style
#      (green fill / green font).
# ---------------------------------------------------------------------
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$ws.Range("B5:B11").PasteSpecial(-4122) | Out-Null
$ws.Range("B13:B17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("F2").Style = "Good"
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F3:F17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Rows that picked up the larger 12pt font now need the taller row
# height that the already-formatted rows (2-4, 12) were using.
$ws.Rows("5:11").RowHeight = 16.5
$ws.Rows("13:17").RowHeight = 16.5

# ---------------------------------------------------------------------
# 4. Move the active selection to B24.
# ---------------------------------------------------------------------
$ws.Range("B24").Select() | Out-Null
